$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.659.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.594.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0836"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.582.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.635.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.692"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.295.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.792"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.729.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.902"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0503"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.80%  "
